# The deck ships two theme parts:
#   ppt/theme/theme1.xml  - "Integral"     (used by the slide master -> every slide)
#   ppt/theme/theme2.xml  - "Office Theme" (used by the notes master)
#
# The target commit swaps the two themes' contents: theme1.xml becomes the
# "Office Theme" colour scheme and theme2.xml becomes "Integral". The font
# scheme and format scheme are identical between the two themes, so the
# whole visible effect of the swap is a change of the 12 theme colours (and
# the cosmetic theme/colour-scheme name) used by the slide master / slides.
#
# Apply the new ("Office Theme") colours to the presentation's theme colour
# scheme via the Slide.ThemeColorScheme COM surface, which is backed by
# ppt/theme/theme1.xml.

function Hex-ToRgbValue($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# ThemeColorScheme item order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1,
# 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$tcs.Item(1).RGB  = Hex-ToRgbValue "000000"   # dk1
$tcs.Item(2).RGB  = Hex-ToRgbValue "FFFFFF"   # lt1
$tcs.Item(3).RGB  = Hex-ToRgbValue "44546A"   # dk2
$tcs.Item(4).RGB  = Hex-ToRgbValue "E7E6E6"   # lt2
$tcs.Item(5).RGB  = Hex-ToRgbValue "5B9BD5"   # accent1
$tcs.Item(6).RGB  = Hex-ToRgbValue "ED7D31"   # accent2
$tcs.Item(7).RGB  = Hex-ToRgbValue "A5A5A5"   # accent3
$tcs.Item(8).RGB  = Hex-ToRgbValue "FFC000"   # accent4
$tcs.Item(9).RGB  = Hex-ToRgbValue "4472C4"   # accent5
$tcs.Item(10).RGB = Hex-ToRgbValue "70AD47"   # accent6
$tcs.Item(11).RGB = Hex-ToRgbValue "0563C1"   # hlink
$tcs.Item(12).RGB = Hex-ToRgbValue "954F72"   # folHlink

# Cosmetic: rename the colour/theme scheme to match the new "Office" theme.
try { $tcs.Name = "Office" } catch {}
